$p = $ppt.ActivePresentation
$p.EmbedTrueTypeFonts = $true
